# Rename "HKID" label to "HKID OR PASSPORT" on the Info sheet, and remove
# the now-redundant "RESIDENTIAL NUMBER" / "OFFICE NUMBER" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

$ws.Range("A9").Value = "HKID OR PASSPORT"

# Delete the "RESIDENTIAL NUMBER" (row 25) and "OFFICE NUMBER" (row 26) rows,
# shifting everything below them up.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(25).Delete()
